# Auto-generated update of cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$style = $cell.Style
$cell.Value = "'42.229.30"
$cell.Style = $style

$cell = $ws.Range("E2")
$style = $cell.Style
$cell.Value = "'  +0.25%  "
$cell.Style = $style

$cell = $ws.Range("D3")
$style = $cell.Style
$cell.Value = "'2.268.98"
$cell.Style = $style

$cell = $ws.Range("E3")
$style = $cell.Style
$cell.Value = "'  -0.45%  "
$cell.Style = $style

$cell = $ws.Range("E4")
$style = $cell.Style
$cell.Value = "'  -0.03%  "
$cell.Style = $style

$cell = $ws.Range("D5")
$style = $cell.Style
$cell.Value = "'307.03"
$cell.Style = $style

$cell = $ws.Range("E5")
$style = $cell.Style
$cell.Value = "'  +0.65%  "
$cell.Style = $style

$cell = $ws.Range("D6")
$style = $cell.Style
$cell.Value = "'97.02"
$cell.Style = $style

$cell = $ws.Range("E6")
$style = $cell.Style
$cell.Value = "'  +3.00%  "
$cell.Style = $style

$cell = $ws.Range("D7")
$style = $cell.Style
$cell.Value = "'0.527"
$cell.Style = $style

$cell = $ws.Range("E7")
$style = $cell.Style
$cell.Value = "'  -0.84%  "
$cell.Style = $style

$cell = $ws.Range("E8")
$style = $cell.Style
$cell.Value = "'  +0.01%  "
$cell.Style = $style

$cell = $ws.Range("E9")
$style = $cell.Style
$cell.Value = "'  +0.83%  "
$cell.Style = $style

$cell = $ws.Range("D10")
$style = $cell.Style
$cell.Value = "'35.32"
$cell.Style = $style

$cell = $ws.Range("E10")
$style = $cell.Style
$cell.Value = "'  +3.62%  "
$cell.Style = $style

$cell = $ws.Range("D11")
$style = $cell.Style
$cell.Value = "'0.0791"
$cell.Style = $style

$cell = $ws.Range("E11")
$style = $cell.Style
$cell.Value = "'  -1.53%  "
$cell.Style = $style

$cell = $ws.Range("E12")
$style = $cell.Style
$cell.Value = "'  -0.02%  "
$cell.Style = $style

$cell = $ws.Range("D13")
$style = $cell.Style
$cell.Value = "'6.88"
$cell.Style = $style

$cell = $ws.Range("E13")
$style = $cell.Style
$cell.Value = "'  +3.14%  "
$cell.Style = $style

$cell = $ws.Range("D14")
$style = $cell.Style
$cell.Value = "'2.620.50"
$cell.Style = $style

$cell = $ws.Range("E14")
$style = $cell.Style
$cell.Value = "'  -0.45%  "
$cell.Style = $style

$cell = $ws.Range("D15")
$style = $cell.Style
$cell.Value = "'14.78"
$cell.Style = $style

$cell = $ws.Range("E15")
$style = $cell.Style
$cell.Value = "'  +2.81%  "
$cell.Style = $style

$cell = $ws.Range("D16")
$style = $cell.Style
$cell.Value = "'2.263.15"
$cell.Style = $style

$cell = $ws.Range("E16")
$style = $cell.Style
$cell.Value = "'  -0.71%  "
$cell.Style = $style

$cell = $ws.Range("E17")
$style = $cell.Style
$cell.Value = "'  +0.20%  "
$cell.Style = $style

$cell = $ws.Range("D18")
$style = $cell.Style
$cell.Value = "'42.097.45"
$cell.Style = $style

$cell = $ws.Range("E18")
$style = $cell.Style
$cell.Value = "'  +0.15%  "
$cell.Style = $style

$cell = $ws.Range("D19")
$style = $cell.Style
$cell.Value = "'12.43"
$cell.Style = $style

$cell = $ws.Range("E19")
$style = $cell.Style
$cell.Value = "'  -2.56%  "
$cell.Style = $style

$cell = $ws.Range("E20")
$style = $cell.Style
$cell.Value = "'  -1.21%  "
$cell.Style = $style

$cell = $ws.Range("D21")
$style = $cell.Style
$cell.Value = "'6.04"
$cell.Style = $style

$cell = $ws.Range("E21")
$style = $cell.Style
$cell.Value = "'  +0.64%  "
$cell.Style = $style

$cell = $ws.Range("D22")
$style = $cell.Style
$cell.Value = "'68.19"
$cell.Style = $style

$cell = $ws.Range("E22")
$style = $cell.Style
$cell.Value = "'  +0.26%  "
$cell.Style = $style

$cell = $ws.Range("D23")
$style = $cell.Style
$cell.Value = "'238.33"
$cell.Style = $style

$cell = $ws.Range("E23")
$style = $cell.Style
$cell.Value = "'  -2.24%  "
$cell.Style = $style

$cell = $ws.Range("D24")
$style = $cell.Style
$cell.Value = "'2.58"
$cell.Style = $style

$cell = $ws.Range("E24")
$style = $cell.Style
$cell.Value = "'  -0.71%  "
$cell.Style = $style

$cell = $ws.Range("E25")
$style = $cell.Style
$cell.Value = "'  +0.16%  "
$cell.Style = $style

$cell = $ws.Range("D26")
$style = $cell.Style
$cell.Value = "'1.00"
$cell.Style = $style

$cell = $ws.Range("D27")
$style = $cell.Style
$cell.Value = "'23.62"
$cell.Style = $style

$cell = $ws.Range("E27")
$style = $cell.Style
$cell.Value = "'  -2.03%  "
$cell.Style = $style

$cell = $ws.Range("D28")
$style = $cell.Style
$cell.Value = "'37.79"
$cell.Style = $style

$cell = $ws.Range("E28")
$style = $cell.Style
$cell.Value = "'  +5.21%  "
$cell.Style = $style

$cell = $ws.Range("E29")
$style = $cell.Style
$cell.Value = "'  -1.90%  "
$cell.Style = $style

$cell = $ws.Range("E30")
$style = $cell.Style
$cell.Value = "'  +0.92%  "
$cell.Style = $style

$cell = $ws.Range("D31")
$style = $cell.Style
$cell.Value = "'161.91"
$cell.Style = $style

$cell = $ws.Range("E31")
$style = $cell.Style
$cell.Value = "'  +0.60%  "
$cell.Style = $style

$cell = $ws.Range("E32")
$style = $cell.Style
$cell.Value = "'  -1.90%  "
$cell.Style = $style

$cell = $ws.Range("E33")
$style = $cell.Style
$cell.Value = "'  +0.05%  "
$cell.Style = $style

$cell = $ws.Range("E34")
$style = $cell.Style
$cell.Value = "'  +3.25%  "
$cell.Style = $style

$cell = $ws.Range("E35")
$style = $cell.Style
$cell.Value = "'  -2.07%  "
$cell.Style = $style

$cell = $ws.Range("D36")
$style = $cell.Style
$cell.Value = "'17.22"
$cell.Style = $style

$cell = $ws.Range("E36")
$style = $cell.Style
$cell.Value = "'  +1.19%  "
$cell.Style = $style

$cell = $ws.Range("E37")
$style = $cell.Style
$cell.Value = "'  -0.21%  "
$cell.Style = $style

$cell = $ws.Range("E38")
$style = $cell.Style
$cell.Value = "'  -3.52%  "
$cell.Style = $style

$cell = $ws.Range("E39")
$style = $cell.Style
$cell.Value = "'  +0.94%  "
$cell.Style = $style

$cell = $ws.Range("E40")
$style = $cell.Style
$cell.Value = "'  -1.61%  "
$cell.Style = $style

$cell = $ws.Range("E41")
$style = $cell.Style
$cell.Value = "'  -4.24%  "
$cell.Style = $style

$cell = $ws.Range("E42")
$style = $cell.Style
$cell.Value = "'  +2.02%  "
$cell.Style = $style

$cell = $ws.Range("D43")
$style = $cell.Style
$cell.Value = "'1.948.26"
$cell.Style = $style

$cell = $ws.Range("E43")
$style = $cell.Style
$cell.Value = "'  -3.72%  "
$cell.Style = $style

$cell = $ws.Range("D44")
$style = $cell.Style
$cell.Value = "'18.91"
$cell.Style = $style

$cell = $ws.Range("E44")
$style = $cell.Style
$cell.Value = "'  -3.65%  "
$cell.Style = $style

$cell = $ws.Range("E45")
$style = $cell.Style
$cell.Value = "'  -0.82%  "
$cell.Style = $style

$cell = $ws.Range("E46")
$style = $cell.Style
$cell.Value = "'  -2.79%  "
$cell.Style = $style

$cell = $ws.Range("E47")
$style = $cell.Style
$cell.Value = "'  -0.82%  "
$cell.Style = $style

$cell = $ws.Range("D48")
$style = $cell.Style
$cell.Value = "'53.66"
$cell.Style = $style

$cell = $ws.Range("E48")
$style = $cell.Style
$cell.Value = "'  +0.27%  "
$cell.Style = $style

$cell = $ws.Range("E49")
$style = $cell.Style
$cell.Value = "'  -0.39%  "
$cell.Style = $style

$cell = $ws.Range("D50")
$style = $cell.Style
$cell.Value = "'92.02"
$cell.Style = $style

$cell = $ws.Range("E50")
$style = $cell.Style
$cell.Value = "'  +0.05%  "
$cell.Style = $style

$cell = $ws.Range("E51")
$style = $cell.Style
$cell.Value = "'  -1.68%  "
$cell.Style = $style
